$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value2 = 2000
$ws.Range("I13").Value2 = 2000
$ws.Range("J13").Value2 = 0
$ws.Range("K13").Value2 = 2000
$ws.Range("L13").Value2 = 0
$ws.Range("M13").Value2 = -1831
$ws.Range("N13").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value2 = 1956773.9
$ws.Range("J17").Value2 = 1956773.9
$ws.Range("L17").Value2 = 5870321.699999999
$ws.Range("N17").Value2 = -5870657.699999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value2 = 1233
$ws.Range("I115").Value2 = 849.5
$ws.Range("J115").Value2 = 2000
$ws.Range("K115").Value2 = 2548.5
$ws.Range("L115").Value2 = 6000
$ws.Range("M115").Value2 = -981.5
$ws.Range("N115").Value2 = -9134

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value2 = 1340
$ws.Range("I137").Value2 = 1109.1666
$ws.Range("J137").Value2 = 2032.5
$ws.Range("K137").Value2 = 3327.4998
$ws.Range("L137").Value2 = 6097.5
$ws.Range("M137").Value2 = -777.4998000000001
$ws.Range("N137").Value2 = -11197.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value2 = 5240.4346
$ws.Range("I45").Value2 = 8486.923000000001
$ws.Range("K45").Value2 = 8486.923000000001
$ws.Range("M45").Value2 = -8109.923000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 2927
$ws.Range("I61").Value2 = 3176.1956
$ws.Range("J61").Value2 = 1653.3334
$ws.Range("K61").Value2 = 3176.1956
$ws.Range("L61").Value2 = 1653.3334
$ws.Range("M61").Value2 = -2964.1956
$ws.Range("N61").Value2 = -2077.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value2 = 1309.8529
$ws.Range("I74").Value2 = 1235.0952
$ws.Range("J74").Value2 = 1430.6154
$ws.Range("K74").Value2 = 1235.0952
$ws.Range("L74").Value2 = 1430.6154
$ws.Range("M74").Value2 = -361.0952
$ws.Range("N74").Value2 = -3178.6154

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value2 = 1309.8529
$ws.Range("I77").Value2 = 1235.0952
$ws.Range("J77").Value2 = 1430.6154
$ws.Range("K77").Value2 = 6175.476
$ws.Range("L77").Value2 = 7153.076999999999
$ws.Range("M77").Value2 = -1807.476
$ws.Range("N77").Value2 = -15889.077

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H119").Value2 = 39333.668
$ws.Range("J119").Value2 = 39333.668
$ws.Range("L119").Value2 = 39333.668
$ws.Range("N119").Value2 = -49009.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value2 = 8557275
$ws.Range("I122").Value2 = 8557275
$ws.Range("J122").Value2 = 0
$ws.Range("K122").Value2 = 25671825
$ws.Range("L122").Value2 = 0
$ws.Range("M122").Value2 = -25669375
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value2 = 58499.332
$ws.Range("J125").Value2 = 58499.332
$ws.Range("L125").Value2 = 58499.332
$ws.Range("N125").Value2 = -68339.33199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value2 = 3229221
$ws.Range("I132").Value2 = 3115.3125
$ws.Range("J132").Value2 = 6670400.5
$ws.Range("K132").Value2 = 9345.9375
$ws.Range("L132").Value2 = 20011201.5
$ws.Range("M132").Value2 = -6815.9375
$ws.Range("N132").Value2 = -20016261.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value2 = 2927
$ws.Range("I136").Value2 = 3176.1956
$ws.Range("J136").Value2 = 1653.3334
$ws.Range("K136").Value2 = 9528.586800000001
$ws.Range("L136").Value2 = 4960.0002
$ws.Range("M136").Value2 = -6978.586800000001
$ws.Range("N136").Value2 = -10060.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value2 = 4295.579
$ws.Range("I134").Value2 = 5006.222
$ws.Range("J134").Value2 = 2551.2727
$ws.Range("K134").Value2 = 15018.666
$ws.Range("L134").Value2 = 7653.8181
$ws.Range("M134").Value2 = -12483.666
$ws.Range("N134").Value2 = -12723.8181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 6429.639
$ws.Range("I31").Value2 = 1780
$ws.Range("J31").Value2 = 15961.4
$ws.Range("K31").Value2 = 1780
$ws.Range("L31").Value2 = 15961.4
$ws.Range("M31").Value2 = -1485
$ws.Range("N31").Value2 = -16551.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value2 = 6429.639
$ws.Range("I34").Value2 = 1780
$ws.Range("J34").Value2 = 15961.4
$ws.Range("K34").Value2 = 1780
$ws.Range("L34").Value2 = 15961.4
$ws.Range("M34").Value2 = -1578
$ws.Range("N34").Value2 = -16365.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value2 = 2006.1777
$ws.Range("I132").Value2 = 1713.8334
$ws.Range("J132").Value2 = 2590.8667
$ws.Range("K132").Value2 = 5141.5002
$ws.Range("L132").Value2 = 7772.6001
$ws.Range("M132").Value2 = -2611.5002
$ws.Range("N132").Value2 = -12832.6001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value2 = 85
$ws.Range("I10").Value2 = 85
$ws.Range("K10").Value2 = 255
$ws.Range("M10").Value2 = -116

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value2 = 42100
$ws.Range("J37").Value2 = 42100
$ws.Range("L37").Value2 = 126300
$ws.Range("N37").Value2 = -126524

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value2 = 1958
$ws.Range("I68").Value2 = 600
$ws.Range("J68").Value2 = 2108.889
$ws.Range("K68").Value2 = 1800
$ws.Range("L68").Value2 = 6326.667
$ws.Range("M68").Value2 = -989
$ws.Range("N68").Value2 = -7948.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value2 = 1958
$ws.Range("I71").Value2 = 600
$ws.Range("J71").Value2 = 2108.889
$ws.Range("K71").Value2 = 5400
$ws.Range("L71").Value2 = 18980.001
$ws.Range("M71").Value2 = -1344
$ws.Range("N71").Value2 = -27092.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value2 = 570.4761999999999
$ws.Range("I107").Value2 = 380
$ws.Range("J107").Value2 = 646.6667
$ws.Range("K107").Value2 = 1140
$ws.Range("L107").Value2 = 1940.0001
$ws.Range("M107").Value2 = 780
$ws.Range("N107").Value2 = -5780.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value2 = 2372.1082
$ws.Range("I132").Value2 = 1985.909
$ws.Range("J132").Value2 = 2938.5334
$ws.Range("K132").Value2 = 5957.727000000001
$ws.Range("L132").Value2 = 8815.600199999999
$ws.Range("M132").Value2 = -3427.727000000001
$ws.Range("N132").Value2 = -13875.6002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value2 = 8266.666999999999
$ws.Range("I14").Value2 = 0
$ws.Range("J14").Value2 = 8266.666999999999
$ws.Range("K14").Value2 = 0
$ws.Range("L14").Value2 = 8266.666999999999
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value2 = -8610.666999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value2 = 1927.5714
$ws.Range("I61").Value2 = 1997.6
$ws.Range("J61").Value2 = 1752.5
$ws.Range("K61").Value2 = 1997.6
$ws.Range("L61").Value2 = 1752.5
$ws.Range("M61").Value2 = -1795.6
$ws.Range("N61").Value2 = -2156.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value2 = 1927.5714
$ws.Range("I113").Value2 = 1997.6
$ws.Range("J113").Value2 = 1752.5
$ws.Range("K113").Value2 = 1997.6
$ws.Range("L113").Value2 = 1752.5
$ws.Range("M113").Value2 = 172.4000000000001
$ws.Range("N113").Value2 = -6092.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value2 = 12147938
$ws.Range("I132").Value2 = 15720084
$ws.Range("J132").Value2 = 2641.9
$ws.Range("K132").Value2 = 47160252
$ws.Range("L132").Value2 = 7925.700000000001
$ws.Range("M132").Value2 = -47157722
$ws.Range("N132").Value2 = -12985.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value2 = 9640.4
$ws.Range("I14").Value2 = 9625
$ws.Range("J14").Value2 = 9647.647000000001
$ws.Range("K14").Value2 = 9625
$ws.Range("L14").Value2 = 9647.647000000001
$ws.Range("M14").Value2 = -9457
$ws.Range("N14").Value2 = -9983.647000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value2 = 1331.2727
$ws.Range("I126").Value2 = 698.5
$ws.Range("J126").Value2 = 1692.8572
$ws.Range("K126").Value2 = 2095.5
$ws.Range("L126").Value2 = 5078.571599999999
$ws.Range("M126").Value2 = 374.5
$ws.Range("N126").Value2 = -10018.5716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value2 = 1190.3684
$ws.Range("I132").Value2 = 913.2083
$ws.Range("J132").Value2 = 1665.5
$ws.Range("K132").Value2 = 2739.6249
$ws.Range("L132").Value2 = 4996.5
$ws.Range("M132").Value2 = -209.6248999999998
$ws.Range("N132").Value2 = -10056.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value2 = 8477012
$ws.Range("I136").Value2 = 2606.9268
$ws.Range("J136").Value2 = 27779822
$ws.Range("K136").Value2 = 7820.780400000001
$ws.Range("L136").Value2 = 83339466
$ws.Range("M136").Value2 = -5270.780400000001
$ws.Range("N136").Value2 = -83344566
